# Applies the edits described by the commit:
# "se agrega hora para guardar archivo de bancos"
#
# Summary of changes to worksheet "Hoja1":
#  - Row 107 gets the data that used to live in row 112 (bank product 9008585507 / BIBO SOLUTIONS SAS)
#  - Row 108 keeps CIMAZ / S.A.S but its "Numero del Producto o Servicio" changes
#  - Row 109 keeps DIANA / CARINA IMPATA RESTREPO but "Codigo del Banco" changes 51 -> 7
#    and its "Numero del Producto o Servicio" changes
#  - Rows 110, 111 and 112 are cleared out (blanked)
#  - The trailing empty rows 232-238 are removed entirely
#  - The view is scrolled up a bit and the selected cell moves from A106 to B106

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 107 ---
$ws.Range("B107").Value = 9008585507
$ws.Range("C107").Value = "BIBO"
$ws.Range("D107").Value = "SOLUTIONS SAS"
$ws.Range("G107").Value = "898098"
$ws.Range("H107").Value = 2359386

# --- Row 108 ---
$ws.Range("G108").Value = "8098098"

# --- Row 109 ---
$ws.Range("E109").Value = 7
$ws.Range("G109").Value = "87897987"

# --- Rows 110-112 are cleared ---
$ws.Range("A110:J112").ClearContents()

# --- Remove trailing empty rows 232-238 ---
$ws.Range("A232:J238").EntireRow.Delete()

# --- Update the view/selection state ---
$ws.Application.ActiveWindow.ScrollRow = 85
$ws.Range("B106").Select()
